$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "94.229.95"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "  -4.16%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.457.84"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "  +2.96%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.999"
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "  -0.17%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "237.36"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "  -6.65%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "639.65"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "  -3.40%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "1.43"
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "  -1.44%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.395"
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "  -8.97%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.999"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "  +0.01%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.959"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "  -7.03%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "3.454.65"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "  +2.88%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "42.11"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "  +0.70%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.197"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "  -5.61%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "6.11"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "  -0.70%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "94.049.99"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "  -4.05%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "4.097.44"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "  +2.39%  "
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "  -3.01%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "8.34"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "  -6.97%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "3.460.54"
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "  +2.55%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "17.67"
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "  -1.76%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "11.33"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "  +4.18%  "
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "  -12.48%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "497.55"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "  -3.82%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "3.23"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "  -5.15%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.0000190"
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "  -6.28%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "6.52"
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "  -2.96%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "90.73"
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "  -7.23%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "11.98"
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = "  -3.68%  "
$ws.Range("B29").Value = "InternetComputer(DFINITY)"
$ws.Range("C29").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "11.71"
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = "  -0.65%  "
$ws.Range("B30").Value = "Dai"
$ws.Range("C30").Value = "https://coinranking.com/coin/MoTuySvg7+dai-dai"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "1.00"
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = "  -0.09%  "
$ws.Range("B31").Value = "PancakeSwap"
$ws.Range("C31").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "2.72"
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = "  +5.33%  "
$ws.Range("B32").Value = "Hedera"
$ws.Range("C32").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.136"
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = "  -8.18%  "
$ws.Range("B33").Value = "Cronos"
$ws.Range("C33").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.180"
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = "  -6.29%  "
$ws.Range("B34").Value = "Binance-PegBSC-USD"
$ws.Range("C34").Value = "https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.00"
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = "  +0.09%  "
$ws.Range("B35").Value = "EthereumClassic"
$ws.Range("C35").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "30.23"
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = "  +4.30%  "
$ws.Range("B36").Value = "PolygonEcosystemToken"
$ws.Range("C36").Value = "https://coinranking.com/coin/iDZ0tG-wI+polygonecosystemtoken-pol"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.556"
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = "  -3.99%  "
$ws.Range("B37").Value = "Bittensor"
$ws.Range("C37").Value = "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "554.50"
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = "  +6.01%  "
$ws.Range("B38").Value = "RenderToken"
$ws.Range("C38").Value = "https://coinranking.com/coin/vfo5XYwcV+rendertoken-render"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "7.61"
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = "  -5.38%  "
$ws.Range("B39").Value = "Fetch.AI"
$ws.Range("C39").Value = "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "1.44"
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "  -5.16%  "
$ws.Range("B40").Value = "USDe"
$ws.Range("C40").Value = "https://coinranking.com/coin/exbfr2U-0+usde-usde"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "1.00"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "  -0.04%  "
$ws.Range("B41").Value = "ARBITRUM"
$ws.Range("C41").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.932"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "  +8.66%  "
$ws.Range("B42").Value = "Kaspa"
$ws.Range("C42").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.150"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "  -2.02%  "
$ws.Range("B43").Value = "WhiteBITCoin"
$ws.Range("C43").Value = "https://coinranking.com/coin/GE4c3_TbB+whitebitcoin-wbt"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "24.05"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "  -1.66%  "
$ws.Range("B44").Value = "ImmutableX"
$ws.Range("C44").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "1.72"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "  -0.29%  "
$ws.Range("B45").Value = "VeChain"
$ws.Range("C45").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.0411"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "  -5.70%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "5.53"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "  -3.51%  "
$ws.Range("B47").Value = "MantraDAO"
$ws.Range("C47").Value = "https://coinranking.com/coin/cTdD8lD-6+mantradao-om"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "3.48"
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "  -4.33%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "2.14"
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "  +3.78%  "
$ws.Range("B49").Value = "dogwifhat"
$ws.Range("C49").Value = "https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "3.19"
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "  -0.21%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "52.89"
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "  -4.48%  "
$ws.Range("B51").Value = "Cosmos"
$ws.Range("C51").Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "8.01"
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "  -8.62%  "
